$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must remain TEXT (they are not
# valid numbers -- e.g. "70.996.29" -- and even the ones that look like
# numbers must keep their exact decimal formatting, e.g. "48.60").
# Force text format before assigning so Excel does not silently convert
# the string to a Number and drop trailing zeros / multi-dot formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.996.29"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.865.13"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "698.07"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.45"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.861.45"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  +5.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.44"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.512.37"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.863.15"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.026.52"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.75"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.25"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.53"
$ws.Range("E22").Value = "  +3.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.09"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.73"
$ws.Range("E26").Value = "  +3.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.33"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.16"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  +8.61%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.68"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.29"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.83"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.31"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.812.70"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.40"
$ws.Range("E39").Value = "  +12.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.08"
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.41"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.04"
$ws.Range("E42").Value = "  +7.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "162.83"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000310"
$ws.Range("E46").Value = "  +5.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.60"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.41"
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.304"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "419.75"
$ws.Range("E50").Value = "  +5.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.69"
$ws.Range("E51").Value = "  +2.16%  "
